$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "Weak" sheet as a copy of "Source" (keeps all the
#    existing formulas / styles / layout), placed right after Source.
# ------------------------------------------------------------------
$src = $wb.Worksheets.Item("Source")
$src.Copy([Type]::Missing, $src)
$weak = $wb.Worksheets.Item("Source (2)")
$weak.Name = "Weak"

# Sheet-scoped defined name "XS" on Weak overriding the workbook-level XS
$weak.Names.Add("XS", "=Weak!`$F`$5")

# Weak uses a different base voltage (F1) than Source
$weak.Range("F1").Value = 211.6

# ------------------------------------------------------------------
# 2. Append the new SC MVA / IBR MVA / SCR / kVs block (rows 6-9) to
#    both "Source" and "Weak". Touch Weak first so the shared-string
#    table records "SC MVA", "IBR MVA", "SCR", "kVs" in that order.
# ------------------------------------------------------------------
foreach ($ws in @($weak, $src)) {
    $ws.Range("E7").Value = "SC MVA"
    $ws.Range("E8").Value = "IBR MVA"
    $ws.Range("E9").Value = "SCR"
    $ws.Range("E6").Value = "kVs"

    $ws.Range("E6:E9").Font.Bold = $true

    $ws.Range("F6").Value = 230
    $ws.Range("F7").Formula = "=F6*F6/F1"
    $ws.Range("F8").Value = 100
    $ws.Range("F9").Formula = "=F7/F8"
}

# ------------------------------------------------------------------
# 3. Selections / active sheet bookkeeping to match the saved file.
# ------------------------------------------------------------------
$src.Range("F9").Select() | Out-Null

$weak.Activate() | Out-Null
$weak.Range("F2").Select() | Out-Null
